$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest values.
# Price cells are forced to Text format so Excel does not reinterpret
# numeric-looking strings (e.g. "242.77") as floating point numbers,
# matching the original inline-string cell type.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.804.55'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.941.21'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.77'
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4880'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2946'
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06882'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('E10').Value = '  +2.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '105.98'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.949.02'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07722'
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.355'
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('E15').Value = '  -2.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '274.33'
$ws.Range('E16').Value = '  -3.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.808.97'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007718'
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.192.34'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.558'
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.723'
$ws.Range('E25').Value = '  -1.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.20'
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.64'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.166'
$ws.Range('E28').Value = '  -1.46%  '
$ws.Range('E30').Value = '  -3.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.567'
$ws.Range('E31').Value = '  -3.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.553'
$ws.Range('E32').Value = '  -2.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.368'
$ws.Range('E33').Value = '  -2.96%  '
$ws.Range('E34').Value = '  -3.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7540'
$ws.Range('E35').Value = '  -1.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.158'
$ws.Range('E36').Value = '  -0.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.000'
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.732'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01995'
$ws.Range('E39').Value = '  -2.62%  '
$ws.Range('E40').Value = '  -2.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.537'
$ws.Range('E41').Value = '  +1.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '77.92'
$ws.Range('E42').Value = '  +7.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.099'
$ws.Range('E43').Value = '  -2.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9051'
$ws.Range('E44').Value = '  +2.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '108.25'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4401'
$ws.Range('E46').Value = '  -1.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9991'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.755'
$ws.Range('E48').Value = '  +3.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.005.58'
$ws.Range('E49').Value = '  +2.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1246'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.275'
$ws.Range('E51').Value = '  -1.36%  '
